$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("URL")

$ws.Range("B15").Value = "https://demoqa.com/buttons"
$ws.Range("A15").Value = "buttonsPage"

$ws.Range("B22").Select()
